$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bitstream-metadata")

# The ROW-ID for this row was pointing at the wrong index (ROW-ID::1);
# fix it to reference the correct row (ROW-ID::2).
$ws.Range("A2").Value = "ROW-ID::2"

# Move the active selection to A3, matching where the cursor lands after
# editing A2 and pressing Enter.
$ws.Range("A3").Select() | Out-Null
